# Extend the "Number of appeals to the Ombudsman" table with three more
# years of data (2021, 2022, 2023) in columns R, S, T, mirroring the
# existing layout/formatting of columns D..Q (one column per year).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (thin bottom-border spacer row under the title) ---------------
# Extend the bordered/formatted row into the new columns (no values, just
# the same look as the existing cells in that row).
$ws.Range("Q2").Copy()
$ws.Range("R2:T2").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 3 (year header row) ----------------------------------------------
$ws.Range("Q3").Copy()
$ws.Range("R3:T3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# --- Row 4 ("Number of written appeals") ----------------------------------
$ws.Range("Q4").Copy()
$ws.Range("R4:T4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620

# --- Row 5 ("Number of positively resolved") ------------------------------
$ws.Range("Q5").Copy()
$ws.Range("R5:T5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264

# Reset the cursor back to the top-left cell (the original file's saved
# selection at F16, well outside the table, is no longer meaningful once
# the sheet only goes to row 5 / column T).
$ws.Range("A1").Select() | Out-Null
